$d = $word.ActiveDocument

# Locate the list item paragraph that reads
# "RightClickable – prefarbovať aj basic plot ITS" and remove it
# entirely (including its paragraph mark), per the commit's removal
# of that TODO entry.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*RightClickable*basic plot ITS*") {
        $p.Range.Delete()
    }
}
